$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.982.63"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.845.42"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.330"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "2.112.31"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.88%  "
$ws.Range("D14").Value = "1.836.41"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").Value = "35.007.81"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "0.0₃0792"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.36%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.124"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.31%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +24.44%  "
$ws.Range("E35").Value = "  +11.08%  "
$ws.Range("E36").Value = "  +9.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.25%  "
$ws.Range("E38").Value = "  +11.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "89.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "1.348.08"
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.04%  "
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0531"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.27%  "
$ws.Range("D48").Value = "2.031.72"
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.64%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -0.02%  "
